# Add the 2020 data column (N) to the insurance-companies summary sheet.
# Mirrors columns D:M (years 2010-2019) which already exist; we extend the
# table by one more year, copying formatting from the adjacent 2019 (M)
# column so the new cells pick up the same styles/number formats/borders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year header (row 3)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2020

# Number of reporting insurance companies (row 4)
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 15

# Insurance premiums, mln soms (row 5)
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 1308.3

# Match the saved selection from the authored workbook.
$ws.Range("N6").Select() | Out-Null
